# Remove the "NO CHANGE" rows (TypeChange column K) that no longer apply.
# Row 8 keeps its real TypeChange value ("VALUESET CHANGES"); every other
# data row (3,4,5,6,7,9,10,11,12) only ever held the placeholder
# "NO CHANGE" value in column K, so clear that cell entirely.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($r in 3,4,5,6,7,9,10,11,12) {
    $ws.Cells.Item($r, 11).ClearContents()
}

# Move the selection to where the edit was made (K9:K12) and scroll the
# viewport over, matching the author's on-screen state after the edit.
$ws.Range("H1").Select()
$ws.Range("K9:K12").Select()
